$wb = $excel.ActiveWorkbook

# --- Sheet "u_MAB" ---
$wsMAB = $wb.Worksheets.Item("u_MAB")
$wsMAB.Range("B15").Value = 0
$wsMAB.Range("A16").Value = 0
$wsMAB.Range("A27").Value = 0.3004191035039643
$wsMAB.Range("B27").Value = 0.1787469213149195
$wsMAB.Range("A40").Value = 0
$wsMAB.Range("B40").Value = 0
$wsMAB.Range("B50").Value = 3.090522235796593
$wsMAB.Range("B51").Value = 0.2742864330750027
$wsMAB.Range("A52").Value = 0.0518270226347731
$wsMAB.Range("B61").Value = 0

# --- Sheet "u_EOH" ---
$wsEOH = $wb.Worksheets.Item("u_EOH")
$wsEOH.Range("A2").Value = -0.2978108300138596
$wsEOH.Range("A3").Value = -0.5771773593421061

# --- Sheet "v_l" ---
$wsVL = $wb.Worksheets.Item("v_l")
$wsVL.Range("A2").Value = 242201.4083119944
$wsVL.Range("A3").Value = 1403128.754780352
$wsVL.Range("A4").Value = 0
